# Implement delayed animations: reorder the per-individual animation columns
# on the "Individuals" sheet so the ANIMATION_ATTACK_* columns immediately
# follow ANIMATION_IDLE, and ANIMATION_HARM/DEATH/CAST move to the end. Also
# set new delay test values for the slash/death columns (rows 2-7), and
# update the active selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")
$ws.Activate()

# --- Row 1 header re-order (columns AT..BA) ---------------------------------
# Before: AT=HARM, AU=DEATH, AV=CAST, AW=SLASH, AX=CHOP, AY=BLUNT, AZ=PIERCE, BA=BOW
# After:  AT=SLASH, AU=CHOP, AV=BLUNT, AW=PIERCE, AX=BOW, AY=HARM, AZ=DEATH, BA=CAST
$ws.Range("AT1").Value = "ANIMATION_ATTACK_SLASH,"
$ws.Range("AU1").Value = "ANIMATION_ATTACK_CHOP,"
$ws.Range("AV1").Value = "ANIMATION_ATTACK_BLUNT,"
$ws.Range("AW1").Value = "ANIMATION_ATTACK_PIERCE,"
$ws.Range("AX1").Value = "ANIMATION_ATTACK_BOW"
$ws.Range("AY1").Value = "ANIMATION_HARM,"
$ws.Range("AZ1").Value = "ANIMATION_DEATH,"
$ws.Range("BA1").Value = "ANIMATION_CAST,"

# --- Data rows 2-7: move values to match the new column order, and set new
#     delay values (5,5,5,5,5,5,-1) for the slash/death animation cells -----
foreach ($row in 2..7) {
    $ws.Range("AT$row").Value = "5,5,5,5,5,5,-1"
    $ws.Range("AU$row").Style = "Normal"
    $ws.Range("AU$row").Value = -1
    $ws.Range("AV$row").Value = -1
    $ws.Range("AW$row").Value = -1
    $ws.Range("AX$row").Value = -1
    $ws.Range("AY$row").Value = -1
    $ws.Range("AZ$row").Value = "5,5,5,5,5,5,-1"
    $ws.Range("BA$row").Value = -1
}

# --- Update the selection left by the editor --------------------------------
$ws.Range("AT7").Select()

Write-Output "done"
